$d = $word.ActiveDocument

# Fix 1: "De Bezoeker moet op de site zijn van park cronestyn " -> split word by
# re-typing the "yn" at the end (simulates the cursor landing there, which is
# also why Word's auto-managed _GoBack bookmark ends up at this spot).
$d.Content.Find.Execute("park cronestyn ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "park cronestyn ", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("cronestyn", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "cronestyn", 2) | Out-Null

# Fix 2: correct the typo "aangemled" -> "aangemeld"
$d.Content.Find.Execute("aangemled", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "aangemeld", 2) | Out-Null
